# Fix upload rapel manfaat
# - corrects "No Transaksi" (column C) values for the 4 data rows
# - corrects "Kode Peserta" (column E) for row 13
# - fills in the "Keterangan" (column N) notes for rows 10-12
# - updates the active selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct "No Transaksi" (col C) values ---
$ws.Range("C10").Value = 566
$ws.Range("C11").Value = 77
$ws.Range("C12").Value = 23
$ws.Range("C13").Value = 521

# --- Correct "Kode Peserta" (col E) for row 13 ---
$ws.Range("E13").Value = 233

# --- Fill "Keterangan" (col N) for rows 10-12 ---
$ws.Range("N10").Value = "tes ubah"
$ws.Range("N11").Value = "tes tes ubah"
$ws.Range("N12").Value = "tes ubah"

# --- Nudge column widths to match the content's new best-fit size ---
$ws.Columns.Item(5).ColumnWidth = 12.333333333333332
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666
$ws.Columns.Item(11).ColumnWidth = 12.833333333333332
$ws.Columns.Item(12).ColumnWidth = 19.333333333333336
$ws.Columns.Item(14).ColumnWidth = 10.166666666666666

# --- Update selection / view position to C13 ---
$ws.Range("C13").Select()
$excel.ActiveWindow.ScrollColumn = 3
